# Refatoracao das metricas de recorrencia trimestral:
# - insere duas novas linhas no topo (2020Q4, 2021Q1)
# - recalcula os valores de todas as linhas existentes
# - adiciona duas novas linhas no final (2025Q1, 2025Q2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere 2 novas linhas logo abaixo do cabecalho, empurrando os dados
# existentes (antigas linhas 2-17) para as linhas 4-19.
$ws.Rows(2).Insert()
$ws.Rows(2).Insert()
$ws.Range("A2:F3").ClearFormats()

# Tabela final completa (trimestre, total_customers, returning_customers,
# new_customers, recurrence_rate) para as linhas 2 a 20.
$data = @(
    @(2,  "2020Q4", 85,  28,  57,  53.84615384615385),
    @(3,  "2021Q1", 195, 47,  148, 55.29411764705883),
    @(4,  "2021Q2", 189, 82,  107, 42.05128205128205),
    @(5,  "2021Q3", 215, 96,  119, 50.79365079365079),
    @(6,  "2021Q4", 217, 129, 88,  60),
    @(7,  "2022Q1", 237, 130, 107, 59.90783410138248),
    @(8,  "2022Q2", 191, 130, 61,  54.85232067510548),
    @(9,  "2022Q3", 177, 126, 51,  65.96858638743456),
    @(10, "2022Q4", 165, 125, 40,  70.62146892655367),
    @(11, "2023Q1", 195, 143, 52,  86.66666666666667),
    @(12, "2023Q2", 201, 153, 48,  78.46153846153847),
    @(13, "2023Q3", 215, 159, 56,  79.1044776119403),
    @(14, "2023Q4", 225, 183, 42,  85.11627906976744),
    @(15, "2024Q1", 255, 186, 69,  82.66666666666667),
    @(16, "2024Q2", 313, 209, 104, 81.96078431372548),
    @(17, "2024Q3", 317, 250, 67,  79.87220447284345),
    @(18, "2024Q4", 319, 255, 64,  80.4416403785489),
    @(19, "2025Q1", 322, 262, 60,  82.13166144200626),
    @(20, "2025Q2", 164, 146, 18,  45.3416149068323)
)

foreach ($row in $data) {
    $r = $row[0]
    $quarter = $row[1]
    $ws.Cells.Item($r, 1).Value = $quarter
    $ws.Cells.Item($r, 2).Value = $quarter
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
